# Edit script: fix manual data, add resource to new attribution entries
$wb = $excel.ActiveWorkbook

# --- Sheet "person matches": trim trailing whitespace from a handful of BDRC-ID cells ---
$ws1 = $wb.Worksheets.Item("person matches")
$trimCells = @("C56", "C57", "C59", "C63", "C73", "C76", "C80", "C83")
foreach ($addr in $trimCells) {
    $cell = $ws1.Range($addr)
    $cur = $cell.Value()
    $cell.Value = $cur.Trim()
}

# --- Sheet "grouped matches": rebuild the BDRC ID / 84000 ID mapping table ---
$ws2 = $wb.Worksheets.Item("grouped matches")

$data = @(
    @('P4CZ16819', '{''eft:sakyaprabha''}'),
    @('P0TMPT007', '{''eft:rnam-par-mi-rtog-pa''}'),
    @('P4CZ15137', '{''eft:kumarakalasa''}'),
    @('P4259', '{''eft:ban-de-dpal-gyi-lhun-po'', ''eft:palgyi-lh-npo'', ''eft:dpal-gyi-lhun-po''}'),
    @('P8266', '{''eft:ch-nyi-tsultrim'', ''eft:dharmatasila''}'),
    @('P8263', '{''eft:leki-d-''}'),
    @('P4255', '{''eft:t-jnanagarbha'', ''eft:yesh-nyingpo'', ''eft:ye-shes-snying-po''}'),
    @('P8220', '{''eft:devacandra''}'),
    @('P00KG07267', '{''eft:sarvajnadeva'', ''eft:sarvanyadeva''}'),
    @('P8182', '{''eft:ska-ba-dpal-brtsegs'', ''eft:paltsek'', ''eft:dpal-brtsegs'', ''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:ban-de-dpal-brtsegs''}'),
    @('P8151', '{''eft:gayadhara''}'),
    @('P8093', '{''eft:kamalagupta''}'),
    @('P8268', '{''eft:buddhaprabha''}'),
    @('P8249', '{''eft:dharmakara''}'),
    @('P4CZ16780', '{''eft:manjusrigarbha''}'),
    @('P4258', '{''eft:dpal-byor''}'),
    @('P8210', '{''eft:danasila''}'),
    @('P8222', '{''eft:jnanasiddhi'', ''eft:jnanasidhi''}'),
    @('P3379', '{''eft:dipamkara-srijnana'', ''eft:dipamkarasrijnana''}'),
    @('P8269', '{''eft:dgon-gling-rma''}'),
    @('P8228', '{''eft:surendrabodhi''}'),
    @('P3214', '{''eft:danasila''}'),
    @('https://lod.dila.edu.tw/resource.php?id=A000089', '{''eft:siladharma''}'),
    @('P8206', '{''eft:celu''}'),
    @('P2548', '{''eft:prajnavarman'', ''eft:prajnavarma''}'),
    @('P0TMP104', '{''eft:punyasambhava''}'),
    @('P0RK8', '{''eft:dharmapala''}'),
    @('P8217', '{''eft:t-jnanagarbha'', ''eft:jnanagarbha''}'),
    @('P8261', '{''eft:munivarma'', ''eft:munivarman''}'),
    @('P8171', '{''eft:dharmasribhadra''}'),
    @('P8209', '{''eft:dzi-na-mi-tra-k-'', ''eft:jinamitra'', ''eft:jinamitra-k-''}'),
    @('P1KG8854', '{''eft:srilendrabodhi'', ''eft:surendrabodhi'', ''eft:silendrabodhi''}'),
    @('P0TMP080', '{''eft:hwa-shang-zab-mo''}'),
    @('P8273', '{''eft:rinchen-tso'', ''eft:rin-chen-tsho''}'),
    @('P3456', '{''eft:tshul-khrims-rgyal-ba''}'),
    @('P4242', '{''eft:sherab-lekpa''}'),
    @('P2637', '{''eft:trakpa-gyaltsen''}'),
    @('P8267', '{''eft:vijayasila''}'),
    @('P8245', '{''eft:buddhakaravarma''}'),
    @('P4263', '{''eft:dge-ba-dpal''}'),
    @('P8183', '{''eft:cog-ro-klu-i-rgyal-mtshan'', ''eft:klu-i-rgyal-mtshan''}'),
    @('P753', '{''eft:rin-chen-bzang-po''}'),
    @('P3285', '{''eft:sakya-yesh-''}'),
    @('P8219', '{''eft:visuddhasimha''}'),
    @('?', '{''eft:sakyasena''}'),
    @('P0TMP092', '{''eft:anandasri-s-''}'),
    @('P5651', '{''eft:pa-tshab-nyi-ma-grags''}'),
    @('P8260', '{''eft:dpal-dbyangs''}'),
    @('P3709', '{''eft:phakpa-sherab''}'),
    @('P8205', '{''eft:zhang-yesh-d-'', ''eft:yesh-d-ye-shes-sde-'', ''eft:band-yesh-de'', ''eft:ye-shes-sde'', ''eft:band-yesh-d-'', ''eft:yesh-d-''}'),
    @('P0TMP098', '{''eft:jinavara''}'),
    @('P2956', '{''eft:krsnapandita''}'),
    @('P8211', '{''eft:vidyakaraprabha''}'),
    @('P8213', '{''eft:t-vidyakarasimha'', ''eft:vidyakarasimha''}'),
    @('P8265', '{''eft:ratnaraksita''}')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 2).Value = $data[$i][0]
    $ws2.Cells.Item($r, 3).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (57-60) left over from the old, longer table
$ws2.Range("A57:C60").EntireRow.Delete() | Out-Null
